$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (strikeout) values for rows 2-70, column G, per regenerated save_data
$kValues = @{
    2 = 0
    3 = 3
    4 = 1
    5 = 1
    6 = 2
    7 = 2
    8 = 0
    9 = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 2
    15 = 2
    16 = 1
    17 = 2
    18 = 2
    19 = 0
    20 = 3
    21 = 3
    22 = 2
    23 = 0
    24 = 3
    25 = 0
    26 = 0
    27 = 3
    28 = 2
    29 = 2
    30 = 1
    31 = 2
    32 = 2
    33 = 1
    34 = 0
    35 = 1
    36 = 2
    37 = 2
    38 = 0
    39 = 2
    40 = 1
    41 = 2
    42 = 1
    43 = 0
    44 = 1
    45 = 0
    46 = 0
    47 = 1
    48 = 1
    49 = 1
    50 = 0
    51 = 2
    52 = 2
    53 = 0
    54 = 0
    55 = 0
    56 = 2
    57 = 2
    58 = 0
    59 = 0
    60 = 3
    61 = 2
    62 = 2
    63 = 0
    64 = 2
    65 = 2
    66 = 1
    67 = 1
    68 = 1
    69 = 3
    70 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
